$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217 (shifts existing rows 217-271 down to 218-272)
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly price record
$ws.Cells.Item(217, 1).Value  = 10
$ws.Cells.Item(217, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(217, 3).Value  = "La Araucanía"
$ws.Cells.Item(217, 4).Value  = 44627
$ws.Cells.Item(217, 5).Value  = 9
$ws.Cells.Item(217, 6).Value  = "Fruta"
$ws.Cells.Item(217, 7).Value  = 100103
$ws.Cells.Item(217, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(217, 9).Value  = 100103004
$ws.Cells.Item(217, 10).Value = "Durazno"
$ws.Cells.Item(217, 11).Value = "Elegant Lady"
$ws.Cells.Item(217, 12).Value = "Primera"
$ws.Cells.Item(217, 13).Value = 100
$ws.Cells.Item(217, 14).Value = 18000
$ws.Cells.Item(217, 15).Value = 18000
$ws.Cells.Item(217, 16).Value = 18000
$ws.Cells.Item(217, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(217, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(217, 19).Value = 1000
$ws.Cells.Item(217, 20).Value = 18
